$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 1060

$ws.Range("C4").Value = 750

$ws.Range("C5").Value = 483

$ws.Range("A6").Value = 'Regione Lombardia'
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = '80050050154'
$ws.Range("C6").Value = 418

$ws.Range("A7").Value = 'Regione del Veneto'
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = '80007580279'
$ws.Range("C7").Value = 409

$ws.Range("A8").Value = 'CREDEMTEL SpA'
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = '01378570350'
$ws.Range("C8").Value = 403

$ws.Range("A9").Value = 'Maggioli SPA'
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = '06188330150'
$ws.Range("C9").Value = 378

$ws.Range("A10").Value = 'Lepida Spa'
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = '02770891204'
$ws.Range("C10").Value = 373

$ws.Range("C12").Value = 320

$ws.Range("C13").Value = 233

$ws.Range("C15").Value = 173

$ws.Range("C16").Value = 150

$ws.Range("C17").Value = 142

$ws.Range("C18").Value = 141

$ws.Range("C19").Value = 131

$ws.Range("C20").Value = 118

$ws.Range("C22").Value = 102

$ws.Range("A23").Value = 'Regione Piemonte'
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = '80087670016'
$ws.Range("C23").Value = 100

$ws.Range("A24").Value = 'ANCITEL'
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = '07196850585'

$ws.Range("C25").Value = 77

$ws.Range("C26").Value = 76

$ws.Range("C27").Value = 60

$ws.Range("C32").Value = 45

$ws.Range("C33").Value = 42

$ws.Range("A36").Value = 'PMPay s.r.l.'
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = '08747230962'
$ws.Range("C36").Value = 38

$ws.Range("A37").Value = 'ROMA CAPITALE'
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = '02438750586'
$ws.Range("C37").Value = 38

$ws.Range("A38").Value = 'SI.net Servizi Informatici S.r.L.'
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = '02743730125'
$ws.Range("C38").Value = 35

$ws.Range("C41").Value = 24

$ws.Range("C42").Value = 24

$ws.Range("A43").Value = 'Nexi SpA'
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = '13212880150'

$ws.Range("A44").Value = 'ANDREANI TRIBUTI srl'
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = '01412920439'

$ws.Range("A45").Value = 'Regione Lazio'
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = '80143490581'

$ws.Range("A46").Value = 'Citta'' Metropolitana di Roma Capitale'
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = '80034390585'
$ws.Range("C46").Value = 18

$ws.Range("A47").Value = 'Comune di Palermo'
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = '80016350821'

$ws.Range("C48").Value = 17

$ws.Range("A51").Value = 'Numera Sistemi e Informatica SpA'
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = '01265230902'
$ws.Range("C51").Value = 12

$ws.Range("A52").Value = 'Si.Form Consulting srl'
$ws.Range("B52").NumberFormat = "@"
$ws.Range("B52").Value = '03943960827'
$ws.Range("C52").Value = 11

$ws.Range("A53").Value = 'Crédit Agricole Group Solutions Società Consortile per azioni'
$ws.Range("B53").NumberFormat = "@"
$ws.Range("B53").Value = '02771790348'
$ws.Range("C53").Value = 11

$ws.Range("A54").Value = 'Aric Agenzia Regionale di Informatica e Committenza'
$ws.Range("B54").NumberFormat = "@"
$ws.Range("B54").Value = '91022630676'
$ws.Range("C54").Value = 10

$ws.Range("A55").Value = 'UBI Banca'
$ws.Range("B55").NumberFormat = "@"
$ws.Range("B55").Value = '03053920165'

$ws.Range("A58").Value = 'Linea Comune Spa'
$ws.Range("B58").NumberFormat = "@"
$ws.Range("B58").Value = '05591950489'

$ws.Range("A59").Value = 'ISWEB S.p.A.'
$ws.Range("B59").NumberFormat = "@"
$ws.Range("B59").Value = '01722270665'

$ws.Range("A60").Value = 'e-SED Società Cooperativa'
$ws.Range("B60").NumberFormat = "@"
$ws.Range("B60").Value = '02695640421'
$ws.Range("C60").Value = 3

$ws.Range("A61").Value = 'ICCREA Banca SpA'
$ws.Range("B61").NumberFormat = "@"
$ws.Range("B61").Value = '04774801007'
$ws.Range("C61").Value = 2

$ws.Range("A62").Value = 'ARCA Servizi s.r.l'
$ws.Range("B62").NumberFormat = "@"
$ws.Range("B62").Value = '09106071005'
$ws.Range("C62").Value = 2

$ws.Range("A63").Value = 'CityPoste Payment Digital S.r.l.'
$ws.Range("B63").NumberFormat = "@"
$ws.Range("B63").Value = '02003750672'
$ws.Range("C63").Value = 2

$ws.Range("A64").Value = 'Ministero dello Sviluppo Economico'
$ws.Range("B64").NumberFormat = "@"
$ws.Range("B64").Value = '80230390587'

$ws.Range("A65").Value = 'Softline srl'
$ws.Range("B65").NumberFormat = "@"
$ws.Range("B65").Value = '12299030150'

$ws.Range("A67").Value = 'Banco BPM Società per Azioni'
$ws.Range("B67").NumberFormat = "@"
$ws.Range("B67").Value = '09722490969'

$ws.Range("A68").Value = 'ARGO SOFTWARE SRL'
$ws.Range("B68").NumberFormat = "@"
$ws.Range("B68").Value = '00838520880'

$ws.Range("A69").Value = 'Engineering Ingegneria Informatica SpA'
$ws.Range("B69").NumberFormat = "@"
$ws.Range("B69").Value = '00967720285'

$ws.Range("A70").Value = 'BANCA MONTE DEI PASCHI DI SIENA'
$ws.Range("B70").NumberFormat = "@"
$ws.Range("B70").Value = '00884060526'

$ws.Range("A71").Value = 'San Marco SPA'
$ws.Range("B71").NumberFormat = "@"
$ws.Range("B71").Value = '04142440728'

$ws.Range("A72").Value = 'Agenzia Italiana del Farmaco - AIFA'
$ws.Range("B72").NumberFormat = "@"
$ws.Range("B72").Value = '97345810580'

$ws.Range("A73").Value = 'MegASP S.r.l.'
$ws.Range("B73").NumberFormat = "@"
$ws.Range("B73").Value = '09898030151'

$ws.Range("A74").Value = 'Società Almaviva S.p.A.'
$ws.Range("B74").NumberFormat = "@"
$ws.Range("B74").Value = '08450891000'
$ws.Range("C74").Value = 1
